$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1741
$ws.Range("F4").Value = 151
$ws.Range("G4").Value = 69
$ws.Range("F5").Value = 403
$ws.Range("F6").Value = 790
$ws.Range("F7").Value = 227
$ws.Range("F8").Value = 1128
$ws.Range("F9").Value = 308
$ws.Range("F13").Value = 179
$ws.Range("F14").Value = 504
$ws.Range("F17").Value = 164
$ws.Range("F18").Value = 2869
$ws.Range("F19").Value = 2603
$ws.Range("F23").Value = 310
$ws.Range("F25").Value = 17
$ws.Range("F26").Value = 5213
$ws.Range("F27").Value = 586
$ws.Range("F29").Value = 17
$ws.Range("F30").Value = 55
$ws.Range("F31").Value = 283
$ws.Range("F32").Value = 1060

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 27
$ws.Range("F4").Value = 1078
$ws.Range("F5").Value = 1078
$ws.Range("F10").Value = 325
$ws.Range("F18").Value = 979
$ws.Range("F27").Value = 3862
$ws.Range("F32").Value = 40

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 76
$ws.Range("F5").Value = 2423
$ws.Range("F6").Value = 1014
$ws.Range("F9").Value = 1275

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2423
$ws.Range("F4").Value = 1741
$ws.Range("F5").Value = 27
$ws.Range("F6").Value = 1014
$ws.Range("F7").Value = 1275
$ws.Range("F10").Value = 151
$ws.Range("G10").Value = 69
$ws.Range("F11").Value = 403
$ws.Range("F12").Value = 790
$ws.Range("F13").Value = 227
$ws.Range("F15").Value = 1128
$ws.Range("F16").Value = 308
$ws.Range("F18").Value = 1078
$ws.Range("F19").Value = 179
$ws.Range("F20").Value = 504
$ws.Range("F22").Value = 164
$ws.Range("F23").Value = 2869
$ws.Range("F24").Value = 2603
$ws.Range("F26").Value = 310
$ws.Range("F27").Value = 325
$ws.Range("F30").Value = 17
$ws.Range("F31").Value = 5213
$ws.Range("F36").Value = 17
$ws.Range("F38").Value = 283
$ws.Range("F44").Value = 1060
$ws.Range("F46").Value = 40
